$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")
$ws.Activate()

# Update the invoice due-date ("fechaVencimiento") values for the two
# inscribed-invoice data rows.
$ws.Range("N2").Value = 65468
$ws.Range("N3").Value = 65468

# Move the view: scroll the window so column V is left-most and select N4
# (matches the new sheetView topLeftCell/selection recorded for "Datos").
$excel.ActiveWindow.ScrollColumn = 22
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N4").Select()
